$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten / rename a handful of menu item titles (new endpoint tests)
$ws.Range("D3").Value = "Сельдь"
$ws.Range("D6").Value = "Горячий "
$ws.Range("D7").Value = "Дайзу "
$ws.Range("D12").Value = "Шем ля Ноблесс"

# Update the selection shown when the sheet is next opened
$ws.Range("A6:H9").Select()
